$wb = $excel.ActiveWorkbook

# Sheet "Add Devices Loop A" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws1.Range("J9").Value = 6
$ws1.Range("J10").Value = 10
$ws1.Activate()
$ws1.Range("A8").Select()

# Sheet "Add_Devices_LoopB_PFI" (sheet2.xml)
$ws2 = $wb.Worksheets.Item("Add_Devices_LoopB_PFI")
$ws2.Range("J9").Value = 10
$ws2.Activate()
$ws2.Range("J9").Select()

# Sheet "Add_Devices_LoopB_FIM" (sheet3.xml)
$ws3 = $wb.Worksheets.Item("Add_Devices_LoopB_FIM")
$ws3.Range("J9").Value = 10
$ws3.Activate()
$ws3.Range("J9").Select()

$ws1.Activate()
